$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.079.58"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.891.27"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.74"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5168"
$ws.Range("E7").Value = "  +2.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3752"
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07212"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.11"
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8998"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07646"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.889.09"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.33"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.228"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008515"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.39"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.133.50"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.053"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.132.81"
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.409"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.286"
$ws.Range("E25").Value = "  +9.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.21"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.733"
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.36"
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.968"
$ws.Range("E30").Value = "  +5.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.775"
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09188"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05045"
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("E34").Value = "  +6.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7743"
$ws.Range("E35").Value = "  +2.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.983"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.279"
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.587"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5603"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01986"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.002"
$ws.Range("E42").Value = "  +5.18%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.641"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.37"
$ws.Range("E44").Value = "  +2.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1512"
$ws.Range("E45").Value = "  +2.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4823"
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.20"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9993"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.593"
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.40"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.85"
$ws.Range("E51").Value = "  +0.91%  "
